$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C407: the page-number lookup came back empty ("NA" removed) ---------
# Matches the existing convention used throughout the sheet for rows with
# no page number (e.g. C2, C3, ...): a blank cell rather than literal text.
$ws.Range("C407").ClearContents()

# --- New scraped rows (408-418) -------------------------------------------
$newRows = @(
    @("2026-01-23", "eaux souterraines", 295, 2),
    @("2026-01-23", "eaux souterraines", 297, 2),
    @("2026-01-23", "eaux de surface",   298, 1),
    @("2026-01-23", "eaux souterraines", 304, 4),
    @("2026-01-23", "eaux souterraines", 305, 1),
    @("2026-01-23", "eaux de surface",   305, 2),
    @("2026-01-23", "ruissellement",     305, 4),
    @("2026-01-23", "eaux souterraines", 311, 1),
    @("2026-01-23", "eaux souterraines", 316, 3),
    @("2026-01-23", "eaux de surface",   317, 1),
    @("2026-01-23", "ruissellement",     317, 1)
)

$startRow = 408
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    # The "Date" column in this sheet stores plain text like "2026-01-23"
    # (not a real Excel date serial). A trailing space keeps Excel's
    # automatic "looks like a date" detection from converting the typed
    # string into a date value while we enter it.
    $ws.Cells.Item($r, 1).Value = ($row[0] + " ")
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
$endRow = $startRow + $newRows.Count - 1

# Strip the trailing space back out of column A via a TRIM() helper column,
# then paste the computed text back as plain values. Copy/PasteSpecial
# values performs a raw copy rather than "typed" user input, so the clean
# "2026-01-23" text lands in A408:A418 without being re-interpreted as a
# date (and without requiring any special text number-format on the cells).
$helper = $ws.Range("F" + $startRow + ":F" + $endRow)
$helper.Formula = "=TRIM(A" + $startRow + ")"
$helper.Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4163)
$helper.ClearContents()
$excel.CutCopyMode = $false
